$wb = $excel.ActiveWorkbook

# The "DeviceList" sheet contains a table of devices in columns B:I (row 1 = device
# names, rows 2-10 = attributes). Column F corresponds to the device
# "SAMSUNG_GalaxyM02_Android_11.0.0_51323" which needs to be removed entirely from
# the device list, shifting every column after it one position to the left.
$ws = $wb.Worksheets.Item("DeviceList")

$ws.Columns.Item(6).EntireColumn.Delete()

# The conditional formatting rules were applied to B2:I2; after the column removal
# they should keep tracking the (now narrower) header row B2:H2.
for ($i = 1; $i -le $ws.Range("B2:H2").FormatConditions.Count; $i++) {
    $fc = $ws.Range("B2:H2").FormatConditions.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("B2:H2"))
}

# Update the view/selection on that sheet to match the new layout.
$ws.Activate()
$ws.Range("E15").Select()
